# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型"
# sheets to reflect a refreshed data scrape (gh-pages output regeneration).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 258
$ws1.Range("F5").Value = 57
$ws1.Range("F8").Value = 54
$ws1.Range("F14").Value = 2023
$ws1.Range("F16").Value = 16
$ws1.Range("F17").Value = 499
$ws1.Range("F18").Value = 468
$ws1.Range("F23").Value = 1522
$ws1.Range("F24").Value = 3428
$ws1.Range("F28").Value = 1113
$ws1.Range("F29").Value = 100
$ws1.Range("F30").Value = 1801
$ws1.Range("F33").Value = 64
$ws1.Range("F37").Value = 649
$ws1.Range("F39").Value = 379

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 258
$ws4.Range("F5").Value = 57
$ws4.Range("F8").Value = 54
$ws4.Range("F14").Value = 2023
$ws4.Range("F17").Value = 16
$ws4.Range("F18").Value = 499
$ws4.Range("F19").Value = 468
$ws4.Range("F24").Value = 1522
$ws4.Range("F25").Value = 3428
$ws4.Range("F29").Value = 1113
$ws4.Range("F30").Value = 100
$ws4.Range("F31").Value = 1801
$ws4.Range("F34").Value = 65
$ws4.Range("F38").Value = 649
$ws4.Range("F40").Value = 379
